$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("B2").Value = "Donnée B2, nouvelle modif 11h40"
$ws.Range("C3").Value = "Push à 11h50"

$ws.Range("D9").Select()
